$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.834.57'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '3.842.76'
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '597.18'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('D7').Value = '3.843.02'
$ws.Range('E7').Value = '  -1.37%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.524'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.165'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('E11').Value = '  -1.32%  '
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '36.77'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = '4.482.81'
$ws.Range('E15').Value = '  -1.50%  '
$ws.Range('D16').Value = '3.819.15'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('D17').Value = '67.850.33'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '18.17'
$ws.Range('E18').Value = '  +7.28%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.43'
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.70'
$ws.Range('E21').Value = '  -4.40%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '466.65'
$ws.Range('E22').Value = '  -3.50%  '
$ws.Range('E23').Value = '  +1.52%  '
$ws.Range('E24').Value = '  -3.46%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '83.04'
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('E26').Value = '  -1.50%  '
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.98'
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('E30').Value = '  -0.21%  '
$ws.Range('D31').Value = '3.987.48'
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('E33').Value = '  -3.34%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '30.88'
$ws.Range('E34').Value = '  -3.07%  '
$ws.Range('D35').Value = '3.811.26'
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('E38').Value = '  -2.66%  '
$ws.Range('E39').Value = '  +0.41%  '
$ws.Range('E40').Value = '  +9.49%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '424.63'
$ws.Range('E43').Value = '  -2.89%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '47.06'
$ws.Range('E46').Value = '  -2.78%  '
$ws.Range('E47').Value = '  +0.94%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '143.36'
$ws.Range('E48').Value = '  +0.92%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '39.59'
$ws.Range('E49').Value = '  +1.51%  '
$ws.Range('E50').Value = '  +10.95%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0355'
$ws.Range('E51').Value = '  +0.26%  '
